$wb = $excel.ActiveWorkbook

function Update-SheetData($ws) {
    $ws.Cells.Item(2,2).Value = "2024-05-04"
    $ws.Cells.Item(2,3).Value = "南昌·Youth动漫美食嘉年华"
    $ws.Cells.Item(2,4).Value = "灌婴路西口朝阳江滩公园内 元亨体育运动中心"
    $ws.Cells.Item(2,5).Value = "2024.05.04 09:00-05.05 20:00"
    $ws.Cells.Item(2,6).Value = 282
    $ws.Cells.Item(2,7).Value = 60
    $ws.Cells.Item(2,8).Value = "https://show.bilibili.com/platform/detail.html?id=84036"
    $ws.Cells.Item(2,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/I5vd5js01712648032400.jpeg"
    $ws.Cells.Item(3,2).Value = "2024-05-12"
    $ws.Cells.Item(3,3).Value = "宜春·BM次元盛典运动番only"
    $ws.Cells.Item(3,4).Value = "鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)"
    $ws.Cells.Item(3,5).Value = "2024.05.12 10:00-05.12 17:00"
    $ws.Cells.Item(3,6).Value = 36
    $ws.Cells.Item(3,7).Value = 55
    $ws.Cells.Item(3,8).Value = "https://show.bilibili.com/platform/detail.html?id=84636"
    $ws.Cells.Item(3,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/sNKPZWMh1713518729449.png"
    $ws.Cells.Item(4,2).Value = "2024-05-18"
    $ws.Cells.Item(4,3).Value = "南昌·花绒万兽首届兽聚"
    $ws.Cells.Item(4,4).Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
    $ws.Cells.Item(4,5).Value = "2024.05.18 09:30-05.19 16:30"
    $ws.Cells.Item(4,6).Value = 137
    $ws.Cells.Item(4,7).Value = 60
    $ws.Cells.Item(4,8).Value = "https://show.bilibili.com/platform/detail.html?id=83689"
    $ws.Cells.Item(4,9).Value = "//i2.hdslb.com/bfs/openplatform/202403/h4iL6IvI1711790121140.jpeg"
    $ws.Cells.Item(5,2).Value = "2024-05-18"
    $ws.Cells.Item(5,3).Value = "赣州·原铁崩only"
    $ws.Cells.Item(5,4).Value = "金岭东大道18号 万达广场(赣州经开店)"
    $ws.Cells.Item(5,5).Value = "2024.05.18 10:00-05.19 17:00"
    $ws.Cells.Item(5,6).Value = 51
    $ws.Cells.Item(5,7).Value = 60
    $ws.Cells.Item(5,8).Value = "https://show.bilibili.com/platform/detail.html?id=84721"
    $ws.Cells.Item(5,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/0n0MQiZh1713505673648.jpeg"
    $ws.Cells.Item(6,2).Value = "2024-05-26"
    $ws.Cells.Item(6,3).Value = "南昌·代号鸢盛花行only"
    $ws.Cells.Item(6,4).Value = "民德路411号 东方豪景花园酒店(民德路店)"
    $ws.Cells.Item(6,5).Value = "2024.05.26 09:30-05.26 17:30"
    $ws.Cells.Item(6,6).Value = 472
    $ws.Cells.Item(6,7).Value = 78
    $ws.Cells.Item(6,8).Value = "https://show.bilibili.com/platform/detail.html?id=82529"
    $ws.Cells.Item(6,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/talOodLW1714030986517.png"
    $ws.Cells.Item(7,2).Value = "2024-06-01"
    $ws.Cells.Item(7,3).Value = "南昌·ACG CLUB动漫游戏嘉年华"
    $ws.Cells.Item(7,4).Value = "火炬五路869号(科技城地铁站3号口步行340米) Ai羽球馆"
    $ws.Cells.Item(7,5).Value = "2024.06.01 10:00-06.01 17:00"
    $ws.Cells.Item(7,6).Value = 1309
    $ws.Cells.Item(7,7).Value = 55
    $ws.Cells.Item(7,8).Value = "https://show.bilibili.com/platform/detail.html?id=84497"
    $ws.Cells.Item(7,9).Value = "//i1.hdslb.com/bfs/openplatform/202404/hZdMDMTZ1713768751631.jpeg"
    $ws.Cells.Item(8,2).Value = "2024-06-08"
    $ws.Cells.Item(8,3).Value = "南昌·CM02动漫游戏博览会"
    $ws.Cells.Item(8,4).Value = "怀玉山大道1315号 南昌绿地国际博览中心"
    $ws.Cells.Item(8,5).Value = "2024.06.08 10:00-06.09 17:00"
    $ws.Cells.Item(8,6).Value = 385
    $ws.Cells.Item(8,7).Value = 65
    $ws.Cells.Item(8,8).Value = "https://show.bilibili.com/platform/detail.html?id=85037"
    $ws.Cells.Item(8,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/qSrEq0US1713947622923.png"
    $ws.Cells.Item(9,2).Value = "2024-06-09"
    $ws.Cells.Item(9,3).Value = "信丰·端午节UPUP动漫展"
    $ws.Cells.Item(9,4).Value = "迎宾大道富华双钻名汇西南侧约200米 诚瑞橙子体育馆"
    $ws.Cells.Item(9,5).Value = "2024.06.09 10:00-06.09 17:00"
    $ws.Cells.Item(9,6).Value = 88
    $ws.Cells.Item(9,7).Value = 48
    $ws.Cells.Item(9,8).Value = "https://show.bilibili.com/platform/detail.html?id=84078"
    $ws.Cells.Item(9,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/Qy0EOl551712651477492.jpeg"
    $ws.Cells.Item(10,2).Value = "2024-06-10"
    $ws.Cells.Item(10,3).Value = "上饶·ETI动漫节"
    $ws.Cells.Item(10,4).Value = "滨江东路与体育馆路交叉口西100米 力加体育综合运动中心"
    $ws.Cells.Item(10,5).Value = "2024.06.10 10:00-06.10 16:00"
    $ws.Cells.Item(10,6).Value = 156
    $ws.Cells.Item(10,7).Value = 55
    $ws.Cells.Item(10,8).Value = "https://show.bilibili.com/platform/detail.html?id=83422"
    $ws.Cells.Item(10,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/N6VdMOuL1713257425864.jpeg"
    $ws.Cells.Item(11,2).Value = "2024-06-10"
    $ws.Cells.Item(11,3).Value = "南昌·LY-COSPLAY大会X运动番PRO2.0（非ONLY）"
    $ws.Cells.Item(11,4).Value = "青山湖南大道260号泰丰轮胎厂进大门走到底左拐 赣A篮球梦时代GANA PARK"
    $ws.Cells.Item(11,5).Value = "2024.06.10 10:00-06.10 17:00"
    $ws.Cells.Item(11,6).Value = 108
    $ws.Cells.Item(11,7).Value = 30
    $ws.Cells.Item(11,8).Value = "https://show.bilibili.com/platform/detail.html?id=84575"
    $ws.Cells.Item(11,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/ScwkijwU1713428452963.jpeg"
    $ws.Cells.Item(12,2).Value = "2024-08-03"
    $ws.Cells.Item(12,3).Value = "南昌·幻梦境国际动漫游戏嘉年华1th"
    $ws.Cells.Item(12,4).Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
    $ws.Cells.Item(12,5).Value = "2024.08.03 09:00-08.04 17:30"
    $ws.Cells.Item(12,6).Value = 156
    $ws.Cells.Item(12,7).Value = 64
    $ws.Cells.Item(12,8).Value = "https://show.bilibili.com/platform/detail.html?id=83980"
    $ws.Cells.Item(12,9).Value = "//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg"
    $ws.Cells.Item(13,2).Value = "2024-08-04"
    $ws.Cells.Item(13,3).Value = "九江·第一届异次元动漫嘉年华"
    $ws.Cells.Item(13,4).Value = "长虹西大道兴城广场99号 九江半岛宾馆"
    $ws.Cells.Item(13,5).Value = "2024.08.04 08:00-08.04 17:00"
    $ws.Cells.Item(13,6).Value = 93
    $ws.Cells.Item(13,7).Value = 35
    $ws.Cells.Item(13,8).Value = "https://show.bilibili.com/platform/detail.html?id=84407"
    $ws.Cells.Item(13,9).Value = "//i1.hdslb.com/bfs/openplatform/202404/e7k26XLV1713262153782.jpeg"
    $ws.Cells.Item(14,2).Value = "2024-08-06"
    $ws.Cells.Item(14,3).Value = "南昌·第一届异次元动漫嘉年华"
    $ws.Cells.Item(14,4).Value = "民德路411号 东方豪景花园酒店(民德路店)"
    $ws.Cells.Item(14,5).Value = "2024.08.06 08:00-08.06 17:00"
    $ws.Cells.Item(14,6).Value = 138
    $ws.Cells.Item(14,7).Value = 40
    $ws.Cells.Item(14,8).Value = "https://show.bilibili.com/platform/detail.html?id=84102"
    $ws.Cells.Item(14,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/0W8gIOYx1712764727806.jpeg"
    $ws.Cells.Item(15,2).Value = "2024-08-08"
    $ws.Cells.Item(15,3).Value = "赣州·第二届异次元动漫嘉年华"
    $ws.Cells.Item(15,4).Value = "金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆"
    $ws.Cells.Item(15,5).Value = "2024.08.08 08:00-08.08 17:00"
    $ws.Cells.Item(15,6).Value = 126
    $ws.Cells.Item(15,7).Value = 45
    $ws.Cells.Item(15,8).Value = "https://show.bilibili.com/platform/detail.html?id=84184"
    $ws.Cells.Item(15,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/F5F9vvqX1712758945373.jpeg"
    $ws.Rows.Item(16).EntireRow.Delete()
    $ws.Rows.Item(16).EntireRow.Delete()
    $ws.Rows.Item(16).EntireRow.Delete()
}

Update-SheetData($wb.Worksheets.Item(1))
Update-SheetData($wb.Worksheets.Item(4))
